# "Generate Report for Archive" — regenerate the localization-status report:
#   1) Status text "Ready for handoff" -> "In Translation" everywhere it appears
#      (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
#   2) The Status column(s) are now narrower (report re-flowed for the new text):
#      Overview!E:F and the Status column (C) on zh-cn / de-de shrink from
#      ~17.22 chars to ~13.41 chars wide.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns (E, F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsOverview.Range("E1:F1").ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
